# "Ez mar a masodik modositas." - second round of edits on this workbook.
#
# 1) Localise a couple of built-in names to Hungarian (best effort - these
#    are read-only/stubbed in some hosts, but harmless to attempt).
# 2) Put 222 into A1.
# 3) Leave the active selection on A2, like the author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hungarian-localise the built-in "Normal" cell style -> "Normal"+accent.
try {
    $normalStyle = $wb.Styles.Item("Normal")
    if ($normalStyle) {
        $normalStyle.Name = [char]0x004E + [char]0x006F + [char]0x0072 + [char]0x006D + [char]0x00E1 + [char]0x006C
    }
} catch {
}

# Hungarian-localise the workbook theme name -> "Office-tema" (with accented e).
try {
    $theme = $wb.Theme
    if ($theme) {
        $theme.Name = [char]0x004F + [char]0x0066 + [char]0x0066 + [char]0x0069 + [char]0x0063 + [char]0x0065 + [char]0x002D + [char]0x0074 + [char]0x00E9 + [char]0x006D + [char]0x0061
    }
} catch {
}

# Hungarian/Japanese-localise the theme's Asian-script (Jpan) font overrides
# on both the major and minor font schemes, best effort. ("Yu Gothic Light"
# -> "You Gothic Light" in Japanese script; "Yu Gothic" -> same w/o " Light".)
try {
    $yuGothicLight = [char]0x6E38 + [char]0x30B4 + [char]0x30B7 + [char]0x30C3 + [char]0x30AF + ' Light'
    $yuGothic = [char]0x6E38 + [char]0x30B4 + [char]0x30B7 + [char]0x30C3 + [char]0x30AF
    $fontScheme = $wb.Theme.ThemeFontScheme
    $fontScheme.MajorFont.EastAsian = $yuGothicLight
    $fontScheme.MinorFont.EastAsian = $yuGothic
} catch {
}

# Write the value that now lives in A1.
$ws.Range("A1").Value = 222

# Match the author's final selection: A2 is the active cell.
$ws.Range("A2").Select() | Out-Null
